# Weekly update: insert 3 new daily-price rows for Plátano (Vega Central
# Mapocho de Santiago) ahead of the existing data block, shifting the prior
# rows 639:747 down to 642:750 (dimension grows from A1:T747 to A1:T750).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down by 3 rows to make room for the new entries.
$ws.Rows("639:641").Insert()

# Values shared by every row in this data block (constant across the sheet).
$commonA = 9
$commonB = "Vega Central Mapocho de Santiago"
$commonC = "Metropolitana"
$commonE = 13
$commonF = "Fruta"
$commonG = 100108
$commonH = "Tropicales y subtropicales"
$commonI = 100108006
$commonJ = "Plátano"
$commonK = "Sin especificar"
$commonQ = "$/caja 20 kilos"
$commonR = "Ecuador"
$commonT = 20

function Set-PlatanoRow {
    param($RowNum, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($RowNum, 1).Value = $commonA
    $ws.Cells.Item($RowNum, 2).Value = $commonB
    $ws.Cells.Item($RowNum, 3).Value = $commonC
    $ws.Cells.Item($RowNum, 4).Value = $Fecha
    $ws.Cells.Item($RowNum, 5).Value = $commonE
    $ws.Cells.Item($RowNum, 6).Value = $commonF
    $ws.Cells.Item($RowNum, 7).Value = $commonG
    $ws.Cells.Item($RowNum, 8).Value = $commonH
    $ws.Cells.Item($RowNum, 9).Value = $commonI
    $ws.Cells.Item($RowNum, 10).Value = $commonJ
    $ws.Cells.Item($RowNum, 11).Value = $commonK
    $ws.Cells.Item($RowNum, 12).Value = $Calidad
    $ws.Cells.Item($RowNum, 13).Value = $Volumen
    $ws.Cells.Item($RowNum, 14).Value = $PrecioMin
    $ws.Cells.Item($RowNum, 15).Value = $PrecioMax
    $ws.Cells.Item($RowNum, 16).Value = $PrecioProm
    $ws.Cells.Item($RowNum, 17).Value = $commonQ
    $ws.Cells.Item($RowNum, 18).Value = $commonR
    $ws.Cells.Item($RowNum, 19).Value = $PrecioKg
    $ws.Cells.Item($RowNum, 20).Value = $commonT
}

Set-PlatanoRow 639 44504 "Pintón"         500  13000 14000 13600 680
Set-PlatanoRow 640 44504 "Primera Maduro" 580  14000 15000 14483 724
Set-PlatanoRow 641 44504 "Primera Pintón" 790  15000 16000 15443 772
